$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Price (D) and Volume 1h % (E) columns for the crypto symbol list,
# refreshed by the scheduled GitHub Actions scraper run.
# Force text formatting before writing so numeric-looking / percent-looking
# strings are stored as literal text, matching the original cell type.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "329.69"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.84%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.17"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.64%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.699"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.16%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08061"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.03%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.026"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "4.15%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.708"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.66%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.522"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.66%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9235"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-2.20%"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.67%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1944"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.43%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.283"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-7.19%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09282"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.37%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.03698"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "5.83%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.1055"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "9.66%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001301"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.83%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006309"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.48%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.381"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.39%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3473"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.50%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.32%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "9.80%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04436"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.55%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001261"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.16%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004364"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.58%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "8.22%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02825"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "16.90%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05460"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.49%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007624"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.07%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009956"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "13.41%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1416"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.60%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002121"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "0.35%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01190"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "9.44%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006713"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.63%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.41%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002999"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-5.39%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002283"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "33.88%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.41%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002004"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.41%"
